$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$rows = @(2,3,4,5,6,9,10,11,12,13,15)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq 79243) {
        $cell.Value2 = 79244
    }
}
